$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q4" right after "总计", which pushes all
#    the other quarter sheets (2022-Q3, 2022-Q2, 2021-Q3, 2021-Q2, 2021-Q1,
#    2020-Q4) one position later without touching their contents.
# ---------------------------------------------------------------------------
$sheetTotal = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheetTotal)
$newSheet.Name = "2022-Q4"

# Header row (same labels used by the other quarterly sheets)
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$newSheet.Range("B1:H1").NumberFormat = "@"
for ($j = 0; $j -lt $headers.Length; $j++) {
    $newSheet.Cells.Item(1, $j + 2).Value = $headers[$j]
}

# Fund holdings data for 2022-Q4 (23 funds)
$fundsData = @(
    @('011201', '财通优势行业轮动混合A', '8.38', '88.19', '7.20', '0.6034', 1),
    @('014029', '浦银安盛红利精选混合C', '9.07', '71.06', '3.71', '0.3365', 4),
    @('519115', '浦银安盛红利精选混合A', '7.60', '71.06', '3.71', '0.2820', 4),
    @('501085', '财通科创主题灵活配置混合（LOF）', '4.19', '89.51', '6.56', '0.2749', 4),
    @('001195', '工银农业产业股票', '5.53', '83.05', '4.70', '0.2599', 2),
    @('519170', '浦银安盛增长动力灵活配置混合A', '6.85', '85.43', '2.24', '0.1534', 8),
    @('519110', '浦银安盛价值成长混合A', '6.74', '90.48', '2.24', '0.1510', 10),
    @('501032', '财通福盛多策略混合（LOF）A', '2.16', '87.95', '6.36', '0.1374', 1),
    @('519125', '浦银安盛消费升级混合A', '2.11', '82.27', '4.77', '0.1006', 6),
    @('540009', '汇丰晋信消费红利股票', '1.72', '90.07', '3.30', '0.0568', 6),
    @('519176', '浦银安盛消费升级混合C', '0.92', '82.27', '4.77', '0.0439', 6),
    @('519120', '浦银安盛新兴产业混合A', '1.72', '91.65', '2.53', '0.0435', 8),
    @('014547', '财通医药鑫选6个月持有期混合A', '0.44', '89.44', '8.22', '0.0362', 1),
    @('501001', '财通多策略精选混合（LOF）', '0.74', '87.17', '4.12', '0.0305', 6),
    @('011202', '财通优势行业轮动混合C', '0.34', '88.19', '7.20', '0.0245', 1),
    @('014548', '财通医药鑫选6个月持有期混合C', '0.28', '89.44', '8.22', '0.0230', 1),
    @('006818', '安信盈利驱动股票A', '0.08', '88.09', '3.20', '0.0026', 10),
    @('006819', '安信盈利驱动股票C', '0.06', '88.09', '3.20', '0.0019', 10),
    @('014003', '浦银安盛增长动力灵活配置混合C', '0.03', '85.43', '2.24', '0.0007', 8),
    @('014011', '浦银安盛价值成长混合C', '0.01', '90.48', '2.24', '0.0002', 10),
    @('014628', '财通福盛多策略混合（LOF）C', '0.00', '87.95', '6.36', 0, 1),
    @('014061', '浦银安盛新兴产业混合C', '0.00', '91.65', '2.53', 0, 8),
    @('960031', '浦银安盛价值成长混合H', '0.00', '90.48', '2.24', 0, 10),
)

$dataRowCount = $fundsData.Length

# Force text formatting on the B:G columns (codes / names / ratios kept as
# text, matching how the other quarterly sheets store them) before writing
# any values, so Excel doesn't silently coerce them to numbers.
$newSheet.Range("B2:G" + ($dataRowCount + 1)).NumberFormat = "@"

for ($i = 0; $i -lt $dataRowCount; $i++) {
    $r = $i + 2
    $row = $fundsData[$i]
    $newSheet.Cells.Item($r, 1).Value = $i
    $newSheet.Cells.Item($r, 2).Value = $row[0]
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $newSheet.Cells.Item($r, 4).Value = $row[2]
    $newSheet.Cells.Item($r, 5).Value = $row[3]
    $newSheet.Cells.Item($r, 6).Value = $row[4]

    # G holds the "持有市值(亿元)" ratio as text, EXCEPT the three smallest
    # holdings whose source value rounded to a bare numeric 0.
    $gVal = $row[5]
    if ($gVal -eq 0) {
        $newSheet.Cells.Item($r, 7).NumberFormat = "General"
    }
    $newSheet.Cells.Item($r, 7).Value = $gVal

    $newSheet.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new row for 2022-Q4 right
#    under the header, shifting the existing quarters down by one row, and
#    renumber the running index in column A.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("总计")
$ws1.Rows.Item(2).Insert()

$ws1.Cells.Item(2, 1).Value = 0
$ws1.Cells.Item(2, 2).Value = "2022-Q4"
$ws1.Cells.Item(2, 3).Value = 23
$ws1.Cells.Item(2, 4).Value = 2.56

for ($r = 3; $r -le 8; $r++) {
    $ws1.Cells.Item($r, 1).Value = $r - 2
}

# Column A carries the bold/centered "index" styling on every data row;
# line it up on the freshly inserted row2 the same way it already is on
# row3 (its index cell came through the insert unstyled).
$ws1.Cells.Item(3, 1).Copy()
$ws1.Cells.Item(2, 1).PasteSpecial(-4122)
